$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at row 17, pushing the existing blank separator row
# (old row 17) and the two rows below it (old 18-19) down to 18-20.
$ws.Range("A17:E17").Insert()

# Copy per-cell formatting from row 16 (same wrap-text / centered style
# used for the new example row) without touching the whole 16384-column
# row the way a full row copy/paste would.
$ws.Range("A16:E16").Copy()
$ws.Range("A17:E17").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Rows.Item(17).RowHeight = 56.25

# Populate the new row with the A2tf example values, in the same order
# the values were originally entered (A, E, B, C, D) so the shared
# string table order matches.
$ws.Range("A17").Value = "char ch[2] // RETURNED"
$ws.Range("E17").Value = "A2tf"
$ws.Range("B17").Value = "[In, Out] char[] signChar,"
$ws.Range("C17").Value = "char[] signChar = {'X'}; [before call]`r`nstring str = new string(signChar); [after]"
$ws.Range("D17").Value = "returns string"

# Match the saved selection state recorded in the workbook.
$ws.Range("D16").Select()
